$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 943
$ws.Range("I111").Value = 899.2857
$ws.Range("J111").Value = 1045
$ws.Range("K111").Value = 2697.8571
$ws.Range("L111").Value = 3135
$ws.Range("M111").Value = 369.1428999999998
$ws.Range("N111").Value = -9269
$ws.Range("H132").Value = 3142.5
$ws.Range("I132").Value = 3076.6667
$ws.Range("J132").Value = 3261
$ws.Range("K132").Value = 9230.000100000001
$ws.Range("L132").Value = 9783
$ws.Range("M132").Value = -6700.000100000001
$ws.Range("N132").Value = -14843
$ws.Range("H138").Value = 2786.6516
$ws.Range("I138").Value = 1578.3636
$ws.Range("J138").Value = 3390.7954
$ws.Range("K138").Value = 4735.0908
$ws.Range("L138").Value = 10172.3862
$ws.Range("M138").Value = 404.9092000000001
$ws.Range("N138").Value = -20452.3862

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 962
$ws.Range("I97").Value = 962
$ws.Range("K97").Value = 962
$ws.Range("M97").Value = -466
$ws.Range("H102").Value = 3699.9285
$ws.Range("I102").Value = 3548.625
$ws.Range("J102").Value = 3901.6667
$ws.Range("K102").Value = 3548.625
$ws.Range("L102").Value = 3901.6667
$ws.Range("M102").Value = -1926.625
$ws.Range("N102").Value = -7145.6667
$ws.Range("H122").Value = 6000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 18000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -22900
$ws.Range("H124").Value = 59860
$ws.Range("J124").Value = 59860
$ws.Range("L124").Value = 59860
$ws.Range("N124").Value = -69680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2155.2222
$ws.Range("J107").Value = 2299.625
$ws.Range("L107").Value = 2299.625
$ws.Range("N107").Value = -6139.625
$ws.Range("H140").Value = 52749.832
$ws.Range("J140").Value = 52749.832
$ws.Range("L140").Value = 52749.832
$ws.Range("N140").Value = -63109.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 875
$ws.Range("J13").Value = 875
$ws.Range("L13").Value = 875
$ws.Range("N13").Value = -1153
$ws.Range("H107").Value = 1140.0938
$ws.Range("I107").Value = 948.5
$ws.Range("K107").Value = 948.5
$ws.Range("M107").Value = 971.5
$ws.Range("H122").Value = 14919.6
$ws.Range("I122").Value = 4866
$ws.Range("K122").Value = 14598
$ws.Range("M122").Value = -12148
$ws.Range("H132").Value = 4728.1333
$ws.Range("I132").Value = 3556.182
$ws.Range("J132").Value = 7951
$ws.Range("K132").Value = 10668.546
$ws.Range("L132").Value = 23853
$ws.Range("M132").Value = -8138.545999999998
$ws.Range("N132").Value = -28913
$ws.Range("H134").Value = 3602.0527
$ws.Range("I134").Value = 2554.4119
$ws.Range("K134").Value = 7663.2357
$ws.Range("M134").Value = -5128.2357

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 19015
$ws.Range("J74").Value = 19015
$ws.Range("L74").Value = 57045
$ws.Range("N74").Value = -59167
$ws.Range("H75").Value = 76927064
$ws.Range("I75").Value = 250002000
$ws.Range("J75").Value = 4868
$ws.Range("K75").Value = 750006000
$ws.Range("L75").Value = 14604
$ws.Range("M75").Value = -750005002
$ws.Range("N75").Value = -16600
$ws.Range("H76").Value = 1967
$ws.Range("I76").Value = 1967
$ws.Range("K76").Value = 5901
$ws.Range("M76").Value = -5518
$ws.Range("H77").Value = 19015
$ws.Range("J77").Value = 19015
$ws.Range("L77").Value = 171135
$ws.Range("N77").Value = -181743
$ws.Range("H78").Value = 76927064
$ws.Range("I78").Value = 250002000
$ws.Range("J78").Value = 4868
$ws.Range("K78").Value = 2250018000
$ws.Range("L78").Value = 43812
$ws.Range("M78").Value = -2250013008
$ws.Range("N78").Value = -53796
$ws.Range("H79").Value = 1967
$ws.Range("I79").Value = 1967
$ws.Range("K79").Value = 5901
$ws.Range("M79").Value = -4575
$ws.Range("H80").Value = 4543.1113
$ws.Range("I80").Value = 1698.3334
$ws.Range("J80").Value = 5965.5
$ws.Range("K80").Value = 5095.0002
$ws.Range("L80").Value = 17896.5
$ws.Range("M80").Value = -4159.0002
$ws.Range("N80").Value = -19768.5
$ws.Range("H82").Value = 3011
$ws.Range("I82").Value = 3011
$ws.Range("K82").Value = 9033
$ws.Range("M82").Value = -8627
$ws.Range("H83").Value = 4543.1113
$ws.Range("I83").Value = 1698.3334
$ws.Range("J83").Value = 5965.5
$ws.Range("K83").Value = 15285.0006
$ws.Range("L83").Value = 53689.5
$ws.Range("M83").Value = -10605.0006
$ws.Range("N83").Value = -63049.5
$ws.Range("H85").Value = 3011
$ws.Range("I85").Value = 3011
$ws.Range("K85").Value = 9033
$ws.Range("M85").Value = -7629
$ws.Range("H98").Value = 1508.45
$ws.Range("J98").Value = 1741.1
$ws.Range("L98").Value = 5223.299999999999
$ws.Range("N98").Value = -8219.299999999999
$ws.Range("H119").Value = 7394.8
$ws.Range("I119").Value = 1861.25
$ws.Range("K119").Value = 5583.75
$ws.Range("M119").Value = -745.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1668.2881
$ws.Range("J102").Value = 3244.5386
$ws.Range("L102").Value = 3244.5386
$ws.Range("N102").Value = -6488.5386
$ws.Range("H118").Value = 23291.5
$ws.Range("J118").Value = 23291.5
$ws.Range("L118").Value = 23291.5
$ws.Range("N118").Value = -26605.5
$ws.Range("H123").Value = 29586.889
$ws.Range("J123").Value = 29586.889
$ws.Range("L123").Value = 29586.889
$ws.Range("N123").Value = -34486.889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 2313
$ws.Range("I32").Value = 2313
$ws.Range("K32").Value = 2313
$ws.Range("M32").Value = -1996
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H93").Value = 1970.75
$ws.Range("I93").Value = 1994.1
$ws.Range("K93").Value = 1994.1
$ws.Range("M93").Value = -746.0999999999999
$ws.Range("H109").Value = 44796.668
$ws.Range("J109").Value = 44796.668
$ws.Range("L109").Value = 44796.668
$ws.Range("N109").Value = -47570.668
$ws.Range("H121").Value = 63157
$ws.Range("J121").Value = 63157
$ws.Range("L121").Value = 63157
$ws.Range("N121").Value = -66651
$ws.Range("H122").Value = 293761.84
$ws.Range("I122").Value = 2002999.5
$ws.Range("K122").Value = 6008998.5
$ws.Range("M122").Value = -6006548.5
$ws.Range("H132").Value = 11579.315
$ws.Range("I132").Value = 10465.583
$ws.Range("J132").Value = 13488.571
$ws.Range("K132").Value = 31396.749
$ws.Range("L132").Value = 40465.713
$ws.Range("M132").Value = -28866.749
$ws.Range("N132").Value = -45525.713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 79106.5
$ws.Range("J138").Value = 79106.5
$ws.Range("L138").Value = 79106.5
$ws.Range("N138").Value = -89386.5
